# Apply "Add data for 2022-07-07" update:
#  - Rename sheet from "Through 2022-06-28" to "Through 2022-06-29"
#  - Update the header label for the rolling current-month column (B1)
#  - Bump several counts in the current-month column (B) and other day-of-week
#    columns by 1 for a handful of neighborhoods, reflecting one more day of
#    incident data being folded into the month-to-date totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet / sheet tab.
$ws.Name = "Through 2022-06-29"

# Update the month-to-date column header text (shared string).
$ws.Range("B1").Value = "June 2022 (through June 29)"

# Row 4 (North Lawndale)
$ws.Range("H4").Value = 7
$ws.Range("N4").Value = 5
$ws.Range("T4").Value = 3

# Row 6 (Humboldt Park)
$ws.Range("H6").Value = 8
$ws.Range("N6").Value = 5
$ws.Range("Z6").Value = 3
$ws.Range("AL6").Value = 3

# Row 8 (Logan Square)
$ws.Range("H8").Value = 2

# Row 10 (Garfield Park)
$ws.Range("Z10").Value = 4

# Row 14 (Austin)
$ws.Range("B14").Value = 11
$ws.Range("H14").Value = 14
$ws.Range("Z14").Value = 5

# Row 17 (West Loop ... )
$ws.Range("B17").Value = 3

# Row 21
$ws.Range("B21").Value = 2

# Row 25 - new value added in previously empty cell
$ws.Range("B25").Value = 1

# Row 30
$ws.Range("T30").Value = 2

# Row 33
$ws.Range("H33").Value = 3

# Row 38
$ws.Range("H38").Value = 2

# Row 41
$ws.Range("N41").Value = 2
$ws.Range("AF41").Value = 4

# Row 59
$ws.Range("AL59").Value = 2

# Row 70
$ws.Range("N70").Value = 2

# Row 71
$ws.Range("B71").Value = 4

# Row 94
$ws.Range("B94").Value = 4

# Row 95
$ws.Range("H95").Value = 2

# Row 96 - new value added in previously empty cell
$ws.Range("T96").Value = 1
